# edit.ps1 - PowerPoint COM-interop script
#
# Reproduces the authored change:
#   1. Three tables (on slides 14, 15 and 16) get their table style
#      switched from {9A7E5B26-370E-4E81-815E-FF3A74368A7B} to
#      {E279C93B-7A8E-4725-AC3F-8D36057C1718}.
#   2. The presentation's main theme (ppt/theme/theme1.xml, currently the
#      "Integral" / "Red Violet" palette used by the Slide Master) is
#      recoloured to the stock Office palette that previously lived only
#      in the Notes Master's theme (ppt/theme/theme2.xml, "Office Theme").
#      Table styles, fonts and effects are already identical between the
#      two theme parts, so swapping the 12 scheme colours reproduces the
#      observable content of the target theme1.xml.

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newStyleId = "{E279C93B-7A8E-4725-AC3F-8D36057C1718}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Theme colours --------------------------------------------------
# Office theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) expressed as OLE BGR-packed RGB() integers, in the order
# expected by ThemeColorScheme.Colors(index).
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
